$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.364.35"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "2.366.68"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  +6.11%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "2.812.68"
$ws.Range("D15").Value = "57.371.69"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.374.10"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "329.45"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.40"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +14.69%  "
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +11.31%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.19"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.26"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.59"
$ws.Range("D32").ClearFormats()
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.917"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.88"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.62"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.38%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.34"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("E50").Value = "  +4.91%  "
$ws.Range("E51").Value = "  +1.55%  "
